# Update "想去人数" (number of people interested) values in column F
# across the four worksheets, per the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(3, 6).Value = 652
$ws1.Cells.Item(4, 6).Value = 816
$ws1.Cells.Item(5, 6).Value = 501
$ws1.Cells.Item(6, 6).Value = 392
$ws1.Cells.Item(7, 6).Value = 479
$ws1.Cells.Item(8, 6).Value = 891
$ws1.Cells.Item(9, 6).Value = 118
$ws1.Cells.Item(10, 6).Value = 834
$ws1.Cells.Item(11, 6).Value = 656
$ws1.Cells.Item(12, 6).Value = 127
$ws1.Cells.Item(13, 6).Value = 42
$ws1.Cells.Item(14, 6).Value = 55
$ws1.Cells.Item(15, 6).Value = 731
$ws1.Cells.Item(16, 6).Value = 219
$ws1.Cells.Item(17, 6).Value = 518
$ws1.Cells.Item(18, 6).Value = 464
$ws1.Cells.Item(19, 6).Value = 1236
$ws1.Cells.Item(21, 6).Value = 959
$ws1.Cells.Item(22, 6).Value = 2704
$ws1.Cells.Item(23, 6).Value = 1169
$ws1.Cells.Item(24, 6).Value = 618
$ws1.Cells.Item(25, 6).Value = 141
$ws1.Cells.Item(26, 6).Value = 1205
$ws1.Cells.Item(28, 6).Value = 899
$ws1.Cells.Item(29, 6).Value = 95
$ws1.Cells.Item(30, 6).Value = 1238

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 489
$ws2.Cells.Item(11, 6).Value = 12
$ws2.Cells.Item(12, 6).Value = 18

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 701

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 701
$ws4.Cells.Item(4, 6).Value = 652
$ws4.Cells.Item(5, 6).Value = 816
$ws4.Cells.Item(6, 6).Value = 501
$ws4.Cells.Item(8, 6).Value = 392
$ws4.Cells.Item(9, 6).Value = 479
$ws4.Cells.Item(10, 6).Value = 489
$ws4.Cells.Item(11, 6).Value = 489
$ws4.Cells.Item(15, 6).Value = 891
$ws4.Cells.Item(16, 6).Value = 118
$ws4.Cells.Item(17, 6).Value = 834
$ws4.Cells.Item(18, 6).Value = 656
$ws4.Cells.Item(19, 6).Value = 127
$ws4.Cells.Item(21, 6).Value = 42
$ws4.Cells.Item(25, 6).Value = 55
$ws4.Cells.Item(26, 6).Value = 12
$ws4.Cells.Item(27, 6).Value = 731
$ws4.Cells.Item(28, 6).Value = 219
$ws4.Cells.Item(29, 6).Value = 518
$ws4.Cells.Item(30, 6).Value = 464
$ws4.Cells.Item(31, 6).Value = 1236
$ws4.Cells.Item(33, 6).Value = 959
$ws4.Cells.Item(34, 6).Value = 2704
$ws4.Cells.Item(35, 6).Value = 1169
$ws4.Cells.Item(36, 6).Value = 618
$ws4.Cells.Item(37, 6).Value = 141
$ws4.Cells.Item(38, 6).Value = 1205
$ws4.Cells.Item(40, 6).Value = 18
$ws4.Cells.Item(41, 6).Value = 899
$ws4.Cells.Item(42, 6).Value = 95
$ws4.Cells.Item(43, 6).Value = 1238

$wb.Save()
